# The deck's sections "SOAP Encoding Rules", "RPC/encoded" and
# "SOAP a iné jazyky" were never populated with slides (the SOAP slides
# were split out into separate presentations) - remove those now-empty
# sections from the section list.

$p  = $ppt.ActivePresentation
$sp = $p.SectionProperties

$namesToRemove = @("SOAP Encoding Rules", "RPC/encoded", "SOAP a iné jazyky")

foreach ($name in $namesToRemove) {
    for ($i = $sp.Count; $i -ge 1; $i--) {
        if ($sp.Name($i) -eq $name) {
            $sp.Delete($i, $false)
            break
        }
    }
}
